{"js": "// Apply each text replacement by searching for the unique old cell text\n// and replacing it with the new text, preserving run formatting.\nconst replacements = [\n  [\"2024-12-22 Sunday\", \"2024-12-23 Monday\"],\n  [\"692\u00f76=115, 2\", \"942\u00f73=314, 0\"],\n  [\"713\u00f72=356, 1\", \"563\u00f74=140, 3\"],\n  [\"831\u00f74=207, 3\", \"896\u00f75=179, 1\"],\n  [\"447\u00f79=49, 6\", \"226\u00f73=75, 1\"],\n  [\"312\u00f77=44, 4\", \"393\u00f74=98, 1\"],\n  [\"273\u00f72=136, 1\", \"267\u00f75=53, 2\"],\n  [\"619\u00f78=77, 3\", \"822\u00f74=205, 2\"],\n  [\"780\u00f72=390, 0\", \"975\u00f77=139, 2\"],\n  [\"528\u00f78=66, 0\", \"868\u00f74=217, 0\"],\n  [\"954\u00f79=106, 0\", \"550\u00f79=61, 1\"],\n  [\"987\u00f77=141, 0\", \"316\u00f77=45, 1\"],\n  [\"365\u00f78=45, 5\", \"842\u00f73=280, 2\"],\n  [\"785\u00f76=130, 5\", \"874\u00f75=174, 4\"],\n  [\"265\u00f77=37, 6\", \"132\u00f75=26, 2\"],\n  [\"904\u00f78=113, 0\", \"619\u00f78=77, 3\"],\n  [\"341\u00f72=170, 1\", \"820\u00f78=102, 4\"],\n  [\"541\u00f78=67, 5\", \"930\u00f73=310, 0\"],\n  [\"249\u00f77=35, 4\", \"531\u00f74=132, 3\"],\n  [\"319\u00f77=45, 4\", \"137\u00f74=34, 1\"],\n  [\"176\u00f74=44, 0\", \"263\u00f72=131, 1\"],\n  [\"108\u00f77=15, 3\", \"306\u00f72=153, 0\"],\n  [\"732\u00f79=81, 3\", \"985\u00f74=246, 1\"],\n  [\"903\u00f75=180, 3\", \"642\u00f72=321, 0\"],\n  [\"159\u00f74=39, 3\", \"401\u00f79=44, 5\"],\n  [\"759\u00f72=379, 1\", \"602\u00f73=200, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$new, 2) | Out-Null\n}\n\nReplace-Text \"2024-12-22 Sunday\" \"2024-12-23 Monday\"\nReplace-Text \"692\u00f76=115, 2\" \"942\u00f73=314, 0\"\nReplace-Text \"713\u00f72=356, 1\" \"563\u00f74=140, 3\"\nReplace-Text \"831\u00f74=207, 3\" \"896\u00f75=179, 1\"\nReplace-Text \"447\u00f79=49, 6\" \"226\u00f73=75, 1\"\nReplace-Text \"312\u00f77=44, 4\" \"393\u00f74=98, 1\"\nReplace-Text \"273\u00f72=136, 1\" \"267\u00f75=53, 2\"\nReplace-Text \"619\u00f78=77, 3\" \"822\u00f74=205, 2\"\nReplace-Text \"780\u00f72=390, 0\" \"975\u00f77=139, 2\"\nReplace-Text \"528\u00f78=66, 0\" \"868\u00f74=217, 0\"\nReplace-Text \"954\u00f79=106, 0\" \"550\u00f79=61, 1\"\nReplace-Text \"987\u00f77=141, 0\" \"316\u00f77=45, 1\"\nReplace-Text \"365\u00f78=45, 5\" \"842\u00f73=280, 2\"\nReplace-Text \"785\u00f76=130, 5\" \"874\u00f75=174, 4\"\nReplace-Text \"265\u00f77=37, 6\" \"132\u00f75=26, 2\"\nReplace-Text \"904\u00f78=113, 0\" \"619\u00f78=77, 3\"\nReplace-Text \"341\u00f72=170, 1\" \"820\u00f78=102, 4\"\nReplace-Text \"541\u00f78=67, 5\" \"930\u00f73=310, 0\"\nReplace-Text \"249\u00f77=35, 4\" \"531\u00f74=132, 3\"\nReplace-Text \"319\u00f77=45, 4\" \"137\u00f74=34, 1\"\nReplace-Text \"176\u00f74=44, 0\" \"263\u00f72=131, 1\"\nReplace-Text \"108\u00f77=15, 3\" \"306\u00f72=153, 0\"\nReplace-Text \"732\u00f79=81, 3\" \"985\u00f74=246, 1\"\nReplace-Text \"903\u00f75=180, 3\" \"642\u00f72=321, 0\"\nReplace-Text \"159\u00f74=39, 3\" \"401\u00f79=44, 5\"\nReplace-Text \"759\u00f72=379, 1\" \"602\u00f73=200, 2\"\n"}
